$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1224:1225").Insert()

$ws.Range("A1224:A1225").Value = 6
$ws.Range("B1224:B1225").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1224:C1225").Value = "Metropolitana"
$ws.Range("D1224").Value = 45166
$ws.Range("D1225").Value = 45166
$ws.Range("E1224:E1225").Value = 13
$ws.Range("F1224:F1225").Value = 100112008
$ws.Range("G1224:G1225").Value = "Coliflor"
$ws.Range("H1224:H1225").Value = "Sin especificar"
$ws.Range("I1224").Value = "Primera"
$ws.Range("I1225").Value = "Segunda"
$ws.Range("J1224").Value = 5000
$ws.Range("K1224").Value = 700
$ws.Range("L1224").Value = 800
$ws.Range("M1224").Value = 734
$ws.Range("J1225").Value = 3600
$ws.Range("K1225").Value = 600
$ws.Range("L1225").Value = 700
$ws.Range("M1225").Value = 647
$ws.Range("N1224:N1225").Value = "$/unidad"
$ws.Range("O1224:O1225").Value = "Región Metropolitana"
$ws.Range("P1224").Value = 734
$ws.Range("P1225").Value = 647
$ws.Range("Q1224:Q1225").Value = 1
$ws.Range("R1224:R1225").Value = "Hortaliza"
